# Fruta / hortaliza, semanal
# Insert a new weekly record as row 241, shifting the existing rows
# (241..335) down by one (to 242..336), and populate the new row with
# the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 241 downward, making room for the new record.
$ws.Rows.Item(241).Insert()

# Fill in the new row with this week's values.
$ws.Cells.Item(241, 1).Value2  = 8
$ws.Cells.Item(241, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(241, 3).Value2  = "Coquimbo"
$ws.Cells.Item(241, 4).Value2  = 44825
$ws.Cells.Item(241, 5).Value2  = 4
$ws.Cells.Item(241, 6).Value2  = 100112003
$ws.Cells.Item(241, 7).Value2  = "Ajo"
$ws.Cells.Item(241, 8).Value2  = "Chino"
$ws.Cells.Item(241, 9).Value2  = "Primera"
$ws.Cells.Item(241, 10).Value2 = 400
$ws.Cells.Item(241, 11).Value2 = 22000
$ws.Cells.Item(241, 12).Value2 = 23000
$ws.Cells.Item(241, 13).Value2 = 22500
$ws.Cells.Item(241, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(241, 15).Value2 = "China"
$ws.Cells.Item(241, 16).Value2 = 2250
$ws.Cells.Item(241, 17).Value2 = 10
$ws.Cells.Item(241, 18).Value2 = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(241, 4).NumberFormat = $ws.Cells.Item(242, 4).NumberFormat
